$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear J3:K12 and J14:K43 (everything except header row 2, row 13, row 44)
$ws.Range("J3:K12").ClearContents()
$ws.Range("J14:K43").ClearContents()

# Update row 13 values
$ws.Range("J13").Value = 1715.970703125
$ws.Range("K13").Value = 120.26168756257545

# Update row 44 values
$ws.Range("J44").Value = 27938.33203125
$ws.Range("K44").Value = 116.16299408456393
